$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows (one at a time, at positions 10 then 12) to make room
# for the new component row and push the totals block (Total/VAT/P&P/Total:)
# down by 2 rows, without leaving a spurious materialised blank row behind.
$ws.Rows("10:10").Insert()
$ws.Rows("12:12").Insert()

# Fill the new component row (row 10) with the stabilising-cap part data.
$ws.Range("A10").Value = "4.7uF ceramic"
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 2.78
$ws.Range("D10").Formula = "=B10*C10"
$ws.Range("E10").Value = "RS"
$ws.Range("F10").Value = "0805"
$ws.Range("G10").Value = "Voltage reg input/output caps. "
$ws.Range("H10").Value = "http://uk.rs-online.com/web/p/ceramic-multilayer-capacitors/7883045/"

# Update the view so the newly added row is visible/selected, matching the
# author's saved cursor position after the edit.
$ws.Range("B1").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("H10").Select()
